$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 17
$ws.Range("B13").Value = 14
$ws.Range("A14").Value = 4
$ws.Range("B14").Value = 8
$ws.Range("A15").Value = 2
$ws.Range("B15").Value = 7

$ws.Range("J16").Select()
